$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert 3 new rows at the end of the synonym table (after current row 112) ---
$ws.Rows("113:115").Insert()

# Copy the formatting (styles) of the last existing synonym row (112) down onto
# the three freshly inserted rows so they match the rest of the table.
$xlPasteFormats = -4122
$ws.Range("A112:B112").Copy() | Out-Null
$ws.Range("A113:B113").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("A114:B114").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("A115:B115").PasteSpecial($xlPasteFormats) | Out-Null
$excel.CutCopyMode = 0

# --- 2. Add the 3 new synonym commands at the bottom of the table ---
$ws.Range("A113").Value = "OPEN"
$ws.Range("B113").Value = "IC"
$ws.Range("A114").Value = "CLOSE"
$ws.Range("B114").Value = "QUIT"
$ws.Range("A115").Value = "LOAD"
$ws.Range("B115").Value = "IC"

# --- 3. Sort the whole synonym table (rows 89-115) alphabetically by column A ---
$xlAscending = 1
$sortRange = $ws.Range("A89:B115")
$keyRange = $ws.Range("A89:A115")
$sortRange.Sort($keyRange, $xlAscending)

# --- 4. Update the view: scroll so row 46 is near the top and select C106 ---
$ws.Range("C106").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 46
$excel.ActiveWindow.ScrollColumn = 1
